$wb = $excel.ActiveWorkbook

# --- 1. Insert new worksheet "2022-Q4" before "2022-Q3" ---
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# NOTE: the sheet that used to be returned by Item(2) is no longer "2022-Q3"
# once the new sheet has been spliced in before it (and the $beforeSheet
# handle itself ends up tracking the newly-inserted sheet) -- re-resolve the
# "2022-Q3" worksheet by name so the Copy() calls below pull from the right
# place.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Copy header row + column-A index style/value from the "2022-Q3" sheet so
# fonts/borders/alignment match the other quarter sheets.
$q3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3.Range("A2").Copy($newSheet.Range("A2:A20"))

$newSheet.Range("B2").Value = "'001173"
$newSheet.Range("C2").Value = "中欧瑾和灵活配置混合 - A"
$newSheet.Range("D2").Value = "'4.56"
$newSheet.Range("E2").Value = "'89.43"
$newSheet.Range("F2").Value = "'4.76"
$newSheet.Range("G2").Value = "'0.2171"
$newSheet.Range("H2").Value = 2
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005994"
$newSheet.Range("C3").Value = "国投瑞银中证500指数量化增强A"
$newSheet.Range("D3").Value = "'13.06"
$newSheet.Range("E3").Value = "'89.93"
$newSheet.Range("F3").Value = "'1.33"
$newSheet.Range("G3").Value = "'0.1737"
$newSheet.Range("H3").Value = 8
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'001174"
$newSheet.Range("C4").Value = "中欧瑾和灵活配置混合 - C"
$newSheet.Range("D4").Value = "'2.63"
$newSheet.Range("E4").Value = "'89.43"
$newSheet.Range("F4").Value = "'4.76"
$newSheet.Range("G4").Value = "'0.1252"
$newSheet.Range("H4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'002906"
$newSheet.Range("C5").Value = "南方中证500量化增强股票A"
$newSheet.Range("D5").Value = "'5.06"
$newSheet.Range("E5").Value = "'91.60"
$newSheet.Range("F5").Value = "'1.45"
$newSheet.Range("G5").Value = "'0.0734"
$newSheet.Range("H5").Value = 4
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'007089"
$newSheet.Range("C6").Value = "国投瑞银中证500指数量化增强C"
$newSheet.Range("D6").Value = "'3.73"
$newSheet.Range("E6").Value = "'89.93"
$newSheet.Range("F6").Value = "'1.33"
$newSheet.Range("G6").Value = "'0.0496"
$newSheet.Range("H6").Value = 8
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'400007"
$newSheet.Range("C7").Value = "东方策略成长混合"
$newSheet.Range("D7").Value = "'1.40"
$newSheet.Range("E7").Value = "'88.54"
$newSheet.Range("F7").Value = "'2.74"
$newSheet.Range("G7").Value = "'0.0384"
$newSheet.Range("H7").Value = 10
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'017288"
$newSheet.Range("C8").Value = "中欧瑾和灵活配置混合 - E"
$newSheet.Range("D8").Value = "'0.45"
$newSheet.Range("E8").Value = "'89.43"
$newSheet.Range("F8").Value = "'4.76"
$newSheet.Range("G8").Value = "'0.0214"
$newSheet.Range("H8").Value = 2
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'002907"
$newSheet.Range("C9").Value = "南方中证500量化增强股票C"
$newSheet.Range("D9").Value = "'1.21"
$newSheet.Range("E9").Value = "'91.60"
$newSheet.Range("F9").Value = "'1.45"
$newSheet.Range("G9").Value = "'0.0175"
$newSheet.Range("H9").Value = 4
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'014588"
$newSheet.Range("C10").Value = "华安中证500指数增强C"
$newSheet.Range("D10").Value = "'0.76"
$newSheet.Range("E10").Value = "'94.54"
$newSheet.Range("F10").Value = "'1.36"
$newSheet.Range("G10").Value = "'0.0103"
$newSheet.Range("H10").Value = 8
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'010154"
$newSheet.Range("C11").Value = "中加中证500指数增强C"
$newSheet.Range("D11").Value = "'0.51"
$newSheet.Range("E11").Value = "'94.31"
$newSheet.Range("F11").Value = "'2.02"
$newSheet.Range("G11").Value = "'0.0103"
$newSheet.Range("H11").Value = 4
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'010153"
$newSheet.Range("C12").Value = "中加中证500指数增强A"
$newSheet.Range("D12").Value = "'0.46"
$newSheet.Range("E12").Value = "'94.31"
$newSheet.Range("F12").Value = "'2.02"
$newSheet.Range("G12").Value = "'0.0093"
$newSheet.Range("H12").Value = 4
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'159620"
$newSheet.Range("C13").Value = "华夏中证智选500成长创新策略ETF"
$newSheet.Range("D13").Value = "'0.32"
$newSheet.Range("E13").Value = "'95.01"
$newSheet.Range("F13").Value = "'1.58"
$newSheet.Range("G13").Value = "'0.0051"
$newSheet.Range("H13").Value = 7
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'014587"
$newSheet.Range("C14").Value = "华安中证500指数增强A"
$newSheet.Range("D14").Value = "'0.33"
$newSheet.Range("E14").Value = "'94.54"
$newSheet.Range("F14").Value = "'1.36"
$newSheet.Range("G14").Value = "'0.0045"
$newSheet.Range("H14").Value = 8
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'013166"
$newSheet.Range("C15").Value = "东兴宸祥量化混合A"
$newSheet.Range("D15").Value = "'0.38"
$newSheet.Range("E15").Value = "'93.88"
$newSheet.Range("F15").Value = "'1.05"
$newSheet.Range("G15").Value = "'0.0040"
$newSheet.Range("H15").Value = 9
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'012297"
$newSheet.Range("C16").Value = "东兴宸瑞量化混合A"
$newSheet.Range("D16").Value = "'0.22"
$newSheet.Range("E16").Value = "'93.65"
$newSheet.Range("F16").Value = "'1.65"
$newSheet.Range("G16").Value = "'0.0036"
$newSheet.Range("H16").Value = 4
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "'005966"
$newSheet.Range("C17").Value = "安信中证500指数增强C"
$newSheet.Range("D17").Value = "'0.16"
$newSheet.Range("E17").Value = "'88.79"
$newSheet.Range("F17").Value = "'0.90"
$newSheet.Range("G17").Value = "'0.0014"
$newSheet.Range("H17").Value = 4
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "'005965"
$newSheet.Range("C18").Value = "安信中证500指数增强A"
$newSheet.Range("D18").Value = "'0.12"
$newSheet.Range("E18").Value = "'88.79"
$newSheet.Range("F18").Value = "'0.90"
$newSheet.Range("G18").Value = "'0.0011"
$newSheet.Range("H18").Value = 4
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "'013167"
$newSheet.Range("C19").Value = "东兴宸祥量化混合C"
$newSheet.Range("D19").Value = "'0.08"
$newSheet.Range("E19").Value = "'93.88"
$newSheet.Range("F19").Value = "'1.05"
$newSheet.Range("G19").Value = "'0.0008"
$newSheet.Range("H19").Value = 9
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "'012298"
$newSheet.Range("C20").Value = "东兴宸瑞量化混合C"
$newSheet.Range("D20").Value = "'0.04"
$newSheet.Range("E20").Value = "'93.65"
$newSheet.Range("F20").Value = "'1.65"
$newSheet.Range("G20").Value = "'0.0007"
$newSheet.Range("H20").Value = 4

# --- 2. Update "总计" summary sheet: insert the 2022-Q4 row at the top of
#        the data (row 2) and push the other quarters down by one row.   ---
$ws1 = $wb.Worksheets.Item("总计")

$ws1.Range("A4:D4").Copy($ws1.Range("A5:D5"))
$ws1.Range("A3:D3").Copy($ws1.Range("A4:D4"))
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 19
$ws1.Range("D2").Value = 0.77

$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
